# Natmi following Dr Hou advice
# Update Lpl-Sdc1 LR-pair edge-weight table: ligand/receptor-expressing cell
# counts changed from 1 to 3 cells, which recomputes the dependent
# expression/specificity/edge-weight statistics for data rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then new values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
# (columns F and L - the detection rates - are unchanged, still 1)
$rows = @(
    @{ Row=2;  E=3; G=71.329076;         H=213.987228;        I=0.3307464087015077; J=0.3307464087015077; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=105.7180836866093;  R=951.4627531794841;  S=0.04636499905115028;  T=0.04636499905115028 },
    @{ Row=3;  E=3; G=71.329076;         H=213.987228;        I=0.3307464087015077; J=0.3307464087015077; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=168.7496384863987;  R=1518.746746377588;  S=0.07400887866542795;  T=0.07400887866542795 },
    @{ Row=4;  E=3; G=71.329076;         H=213.987228;        I=0.3307464087015077; J=0.3307464087015077; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=479.6760765915906;  R=4317.084689324316;  S=0.2103725309849295;  T=0.2103725309849295 },
    @{ Row=5;  E=3; G=129.5615336666667; H=388.684601;        I=0.6007649947142101; J=0.6007649947142101; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=192.0254379677948;  R=1728.228941710153;  S=0.08421699428043308;  T=0.08421699428043307 },
    @{ Row=6;  E=3; G=129.5615336666667; H=388.684601;        I=0.6007649947142101; J=0.6007649947142101; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=306.5154239204413;  R=2758.638815283971;  S=0.1344291046871698;  T=0.1344291046871698 },
    @{ Row=7;  E=3; G=129.5615336666667; H=388.684601;        I=0.6007649947142101; J=0.6007649947142101; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=871.2795907578553;  R=7841.516316820697;  S=0.3821188957466072;  T=0.3821188957466072 },
    @{ Row=8;  E=3; G=14.770314;         H=44.310942;         I=0.06848859658428214; J=0.06848859658428215; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=21.891343321614;    R=197.022089894526;   S=0.009600931807881428; T=0.009600931807881428 },
    @{ Row=9;  E=3; G=14.770314;         H=44.310942;         I=0.06848859658428214; J=0.06848859658428215; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=34.943466081498;    R=314.491194733482;   S=0.01532522833572486;  T=0.01532522833572486 },
    @{ Row=10; E=3; G=14.770314;         H=44.310942;         I=0.06848859658428214; J=0.06848859658428215; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=99.32788516068599;  R=893.950966446174;   S=0.04356243644067585;  T=0.04356243644067587 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("E$n").Value = $r.E
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
}
